$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# TC1 block
$ws.Range("B8").Value = "Lider de Pessoas esta autenticado no sistema e tem permissao para gerenciar Metas de Desempenho"
$ws.Range("B11").Value = "Lider de Pessoas com uma avaliacao selecionada, clica na opcao 'Editar' para modificar a Avaliacao de Desempenho"
$ws.Range("D11").Value = "SYSTEM apresenta o formulario com o campo 'Metas' contendo cada Competencia do perfil avaliado"

# TC2 block (same text repeated)
$ws.Range("B19").Value = "Lider de Pessoas esta autenticado no sistema e tem permissao para gerenciar Metas de Desempenho"
$ws.Range("B22").Value = "Lider de Pessoas com uma avaliacao selecionada, clica na opcao 'Editar' para modificar a Avaliacao de Desempenho"
$ws.Range("D22").Value = "SYSTEM apresenta o formulario com o campo 'Metas' contendo cada Competencia do perfil avaliado"
